$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Calculations")

$newVals = @{
  2 = 161.22999999999999
  3 = 165.19
  4 = 141.15
  5 = 309.69
  6 = 55.06
  7 = 66.59
  8 = 520.71
  9 = 513.66999999999996
  10 = 151.88999999999999
  11 = 150.5
  12 = 219.71
  13 = 306.64
  14 = 526.16999999999996
  15 = 619.85
}

foreach ($row in $newVals.Keys) {
    $ws2.Cells.Item($row, 3).Value = $newVals[$row]
}

$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("E16").NumberFormat = "0"
$ws1.Range("G16:K16").NumberFormat = "0"
$ws1.Range("E17").NumberFormat = "0"
$ws1.Range("G17:K17").NumberFormat = "0"

$ws1.Range("F2").Select()
$ws2.Activate()
$ws2.Range("C2:C15").Select()
